# Updated cryptos list on Tue Oct  3 15:15:12 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain TEXT (mirrors the source file's inlineStr
# cells) without leaving a lasting number-format override on the cell, so
# numeric-looking strings (e.g. "213.71") don't get auto-coerced to numbers.
function Set-TextValue([string]$addr, [string]$value) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "27.359.67"
$ws.Range("E2").Value = "  -3.38%  "

# Row 3 - Ethereum
Set-TextValue "D3" "1.649.20"
$ws.Range("E3").Value = "  -3.85%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.02%  "

# Row 5 - BNB
Set-TextValue "D5" "213.71"
$ws.Range("E5").Value = "  -2.08%  "

# Row 6 - XRP
Set-TextValue "D6" "0.509"
$ws.Range("E6").Value = "  -2.49%  "

# Row 8 - Solana
Set-TextValue "D8" "23.98"
$ws.Range("E8").Value = "  -0.80%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  -1.76%  "

# Row 10 - Dogecoin
Set-TextValue "D10" "0.0614"
$ws.Range("E10").Value = "  -2.81%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  -1.62%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-TextValue "D12" "1.882.92"
$ws.Range("E12").Value = "  -3.89%  "

# Row 13 - WrappedEther
Set-TextValue "D13" "1.656.99"
$ws.Range("E13").Value = "  -3.47%  "

# Row 14 - Polkadot
Set-TextValue "D14" "4.09"
$ws.Range("E14").Value = "  -2.85%  "

# Row 15 - Polygon
Set-TextValue "D15" "0.561"
$ws.Range("E15").Value = "  -0.06%  "

# Row 16 - Litecoin
Set-TextValue "D16" "65.55"
$ws.Range("E16").Value = "  -2.77%  "

# Row 17 - WrappedBTC
Set-TextValue "D17" "27.355.34"

# Row 18 - BitcoinCash
Set-TextValue "D18" "234.50"
$ws.Range("E18").Value = "  -7.44%  "

# Row 19 - ShibaInu
Set-TextValue "D19" "0.0₃0724"
$ws.Range("E19").Value = "  -2.95%  "

# Row 20 - Chainlink
Set-TextValue "D20" "7.49"
$ws.Range("E20").Value = "  -3.65%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  +0.12%  "

# Row 22 - Uniswap
Set-TextValue "D22" "4.40"
$ws.Range("E22").Value = "  -3.68%  "

# Row 23 - Avalanche
Set-TextValue "D23" "9.20"
$ws.Range("E23").Value = "  -4.09%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  -0.21%  "

# Row 25 - Monero
Set-TextValue "D25" "145.77"
$ws.Range("E25").Value = "  -1.48%  "

# Row 26 - Cosmos
Set-TextValue "D26" "7.15"
$ws.Range("E26").Value = "  -3.23%  "

# Row 27 - EthereumClassic
Set-TextValue "D27" "16.02"
$ws.Range("E27").Value = "  -3.51%  "

# Row 28 - BinanceUSD
$ws.Range("E28").Value = "  -0.01%  "

# Row 29 - Stellar
$ws.Range("E29").Value = "  -2.53%  "

# Row 30 - Hedera
Set-TextValue "D30" "0.0496"
$ws.Range("E30").Value = "  -2.93%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  -0.93%  "

# Row 32 - Filecoin
Set-TextValue "D32" "3.29"
$ws.Range("E32").Value = "  -3.12%  "

# Row 33 - Maker
Set-TextValue "D33" "1.454.53"
$ws.Range("E33").Value = "  -1.79%  "

# Row 34 - InternetComputer(DFINITY)
$ws.Range("E34").Value = "  -4.12%  "

# Row 35 - LidoDAOToken
$ws.Range("E35").Value = "  -4.89%  "

# Row 36 - HuobiToken
$ws.Range("E36").Value = "  -0.08%  "

# Row 37 - ARBITRUM
Set-TextValue "D37" "0.909"
$ws.Range("E37").Value = "  -6.24%  "

# Row 38 - ImmutableX
$ws.Range("E38").Value = "  -4.68%  "

# Row 39 - VeChain
$ws.Range("E39").Value = "  -2.97%  "

# Row 40 - WEMIXToken
$ws.Range("E40").Value = "  -0.84%  "

# Row 41 - PaxDollar
$ws.Range("E41").Value = "  +0.14%  "

# Row 42 & 43 - Aave / FraxShare swap places (42 becomes FraxShare, 43 becomes Aave)
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D42" "5.43"
$ws.Range("E42").Value = "  -3.93%  "

$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D43" "65.22"
$ws.Range("E43").Value = "  -6.25%  "

# Row 44 - MXToken
$ws.Range("E44").Value = "  -2.25%  "

# Row 45 - TrustWalletToken
Set-TextValue "D45" "0.786"
$ws.Range("E45").Value = "  -1.92%  "

# Row 46 - RocketPoolETH
Set-TextValue "D46" "1.790.43"
$ws.Range("E46").Value = "  -3.87%  "

# Row 47 - RenderToken
$ws.Range("E47").Value = "  -1.99%  "

# Row 48 - Quant
Set-TextValue "D48" "88.28"
$ws.Range("E48").Value = "  -1.80%  "

# Row 49 - BabyDogeCoin
$ws.Range("E49").Value = "  -4.86%  "

# Row 50 - Algorand
$ws.Range("E50").Value = "  -2.64%  "

# Row 51 - EnergySwap
Set-TextValue "D51" "7.77"
$ws.Range("E51").Value = "  -3.50%  "
